# Fruta / hortaliza, semanal
# Insert two new weekly report rows for "Frutilla" at rows 152-153,
# pushing the existing rows 152-161 down to 154-163.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 152 (shifts 152:161 -> 154:163)
$ws.Rows("152:153").Insert()

# New row 152 data
$ws.Cells.Item(152, 1).Value = 7
$ws.Cells.Item(152, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(152, 3).Value = "Ñuble"
$ws.Cells.Item(152, 4).Value = 44516
$ws.Cells.Item(152, 5).Value = 16
$ws.Cells.Item(152, 6).Value = "Fruta"
$ws.Cells.Item(152, 7).Value = 100101
$ws.Cells.Item(152, 8).Value = "Berries"
$ws.Cells.Item(152, 9).Value = 100112025
$ws.Cells.Item(152, 10).Value = "Frutilla"
$ws.Cells.Item(152, 11).Value = "Sin especificar"
$ws.Cells.Item(152, 12).Value = "Especial"
$ws.Cells.Item(152, 13).Value = 80
$ws.Cells.Item(152, 14).Value = 9000
$ws.Cells.Item(152, 15).Value = 9000
$ws.Cells.Item(152, 16).Value = 9000
$ws.Cells.Item(152, 17).Value = "`$/caja 7 kilos"
$ws.Cells.Item(152, 18).Value = "Provincia de Diguillín"
$ws.Cells.Item(152, 19).Value = 1286
$ws.Cells.Item(152, 20).Value = 7

# New row 153 data
$ws.Cells.Item(153, 1).Value = 7
$ws.Cells.Item(153, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(153, 3).Value = "Ñuble"
$ws.Cells.Item(153, 4).Value = 44516
$ws.Cells.Item(153, 5).Value = 16
$ws.Cells.Item(153, 6).Value = "Fruta"
$ws.Cells.Item(153, 7).Value = 100101
$ws.Cells.Item(153, 8).Value = "Berries"
$ws.Cells.Item(153, 9).Value = 100112025
$ws.Cells.Item(153, 10).Value = "Frutilla"
$ws.Cells.Item(153, 11).Value = "Sin especificar"
$ws.Cells.Item(153, 12).Value = "Primera"
$ws.Cells.Item(153, 13).Value = 160
$ws.Cells.Item(153, 14).Value = 8000
$ws.Cells.Item(153, 15).Value = 8500
$ws.Cells.Item(153, 16).Value = 8250
$ws.Cells.Item(153, 17).Value = "`$/caja 7 kilos"
$ws.Cells.Item(153, 18).Value = "Provincia de Diguillín"
$ws.Cells.Item(153, 19).Value = 1179
$ws.Cells.Item(153, 20).Value = 7

# Make sure the D column (date) of the new rows keeps the same date-time
# number format as the rest of the column.
$ws.Range("D152:D153").NumberFormat = $ws.Range("D151").NumberFormat
